$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 172 <- data formerly on row 174
$ws.Cells.Item(172, 2).Value = 7153759
$ws.Cells.Item(172, 3).Value = 'Azerbaijan Premier League'
$ws.Cells.Item(172, 4).Value = 45430.41666666666
$ws.Cells.Item(172, 5).Value = 'Sabah'
$ws.Cells.Item(172, 6).Value = 'FK Sumqayit'
$ws.Cells.Item(172, 7).Value = 2
$ws.Cells.Item(172, 8).Value = 0
$ws.Cells.Item(172, 9).Value = 1
$ws.Cells.Item(172, 10).Value = 0
$ws.Cells.Item(172, 11).Value = 'H'
$ws.Cells.Item(172, 12).Value = 1.833
$ws.Cells.Item(172, 13).Value = 3.3
$ws.Cells.Item(172, 14).Value = 3.8
$ws.Cells.Item(172, 15).Value = 2.15
$ws.Cells.Item(172, 16).Value = 3.2
$ws.Cells.Item(172, 17).Value = 3.1
$ws.Cells.Item(172, 18).Value = -0.25
$ws.Cells.Item(172, 19).Value = 1.9
$ws.Cells.Item(172, 20).Value = 1.9
$ws.Cells.Item(172, 21).Value = 2.25
$ws.Cells.Item(172, 22).Value = 1.8
$ws.Cells.Item(172, 23).Value = 2
$ws.Cells.Item(172, 24).Value = 1.15
$ws.Cells.Item(172, 25).Value = -1
$ws.Cells.Item(172, 26).Value = -1
$ws.Cells.Item(172, 27).Value = 0.8999999999999999
$ws.Cells.Item(172, 28).Value = -1
$ws.Cells.Item(172, 29).Value = -0.5
$ws.Cells.Item(172, 30).Value = 0.5

# Row 174 <- data formerly on row 172
$ws.Cells.Item(174, 2).Value = 7158118
$ws.Cells.Item(174, 3).Value = 'Azerbaijan Premier League'
$ws.Cells.Item(174, 4).Value = 45430.41666666666
$ws.Cells.Item(174, 5).Value = 'Sabail FC'
$ws.Cells.Item(174, 6).Value = 'Neftchi Baku'
$ws.Cells.Item(174, 7).Value = 0
$ws.Cells.Item(174, 8).Value = 3
$ws.Cells.Item(174, 9).Value = 0
$ws.Cells.Item(174, 10).Value = 0
$ws.Cells.Item(174, 11).Value = 'A'
$ws.Cells.Item(174, 12).Value = 4.5
$ws.Cells.Item(174, 13).Value = 3.7
$ws.Cells.Item(174, 14).Value = 1.6
$ws.Cells.Item(174, 15).Value = 4.2
$ws.Cells.Item(174, 16).Value = 3.5
$ws.Cells.Item(174, 17).Value = 1.7
$ws.Cells.Item(174, 18).Value = 0.75
$ws.Cells.Item(174, 19).Value = 1.9
$ws.Cells.Item(174, 20).Value = 1.9
$ws.Cells.Item(174, 21).Value = 2.75
$ws.Cells.Item(174, 22).Value = 1.825
$ws.Cells.Item(174, 23).Value = 1.975
$ws.Cells.Item(174, 24).Value = -1
$ws.Cells.Item(174, 25).Value = -1
$ws.Cells.Item(174, 26).Value = 0.7
$ws.Cells.Item(174, 27).Value = -1
$ws.Cells.Item(174, 28).Value = 0.8999999999999999
$ws.Cells.Item(174, 29).Value = 0.4125
$ws.Cells.Item(174, 30).Value = -0.5

# Row 176 <- data formerly on row 178
$ws.Cells.Item(176, 2).Value = 7223244
$ws.Cells.Item(176, 3).Value = 'Azerbaijan Premier League'
$ws.Cells.Item(176, 4).Value = 45437.5
$ws.Cells.Item(176, 5).Value = 'Neftchi Baku'
$ws.Cells.Item(176, 6).Value = 'Sabah'
$ws.Cells.Item(176, 7).Value = 0
$ws.Cells.Item(176, 8).Value = 1
$ws.Cells.Item(176, 9).Value = 0
$ws.Cells.Item(176, 10).Value = 0
$ws.Cells.Item(176, 11).Value = 'A'
$ws.Cells.Item(176, 12).Value = 2.1
$ws.Cells.Item(176, 13).Value = 3
$ws.Cells.Item(176, 14).Value = 3.25
$ws.Cells.Item(176, 15).Value = 1.8
$ws.Cells.Item(176, 16).Value = 3.2
$ws.Cells.Item(176, 17).Value = 4
$ws.Cells.Item(176, 18).Value = -0.5
$ws.Cells.Item(176, 19).Value = 1.825
$ws.Cells.Item(176, 20).Value = 1.975
$ws.Cells.Item(176, 21).Value = 2.75
$ws.Cells.Item(176, 22).Value = 1.95
$ws.Cells.Item(176, 23).Value = 1.85
$ws.Cells.Item(176, 24).Value = -1
$ws.Cells.Item(176, 25).Value = -1
$ws.Cells.Item(176, 26).Value = 3
$ws.Cells.Item(176, 27).Value = -1
$ws.Cells.Item(176, 28).Value = 0.9750000000000001
$ws.Cells.Item(176, 29).Value = -1
$ws.Cells.Item(176, 30).Value = 0.8500000000000001

# Row 177 <- data formerly on row 176
$ws.Cells.Item(177, 2).Value = 7217874
$ws.Cells.Item(177, 3).Value = 'Azerbaijan Premier League'
$ws.Cells.Item(177, 4).Value = 45437.5
$ws.Cells.Item(177, 5).Value = 'FK Sumqayit'
$ws.Cells.Item(177, 6).Value = 'FK Gabala'
$ws.Cells.Item(177, 7).Value = 1
$ws.Cells.Item(177, 8).Value = 0
$ws.Cells.Item(177, 9).Value = 0
$ws.Cells.Item(177, 10).Value = 0
$ws.Cells.Item(177, 11).Value = 'H'
$ws.Cells.Item(177, 12).Value = 1.333
$ws.Cells.Item(177, 13).Value = 4.5
$ws.Cells.Item(177, 14).Value = 7
$ws.Cells.Item(177, 15).Value = 1.55
$ws.Cells.Item(177, 16).Value = 3.9
$ws.Cells.Item(177, 17).Value = 4.5
$ws.Cells.Item(177, 18).Value = -1
$ws.Cells.Item(177, 19).Value = 1.975
$ws.Cells.Item(177, 20).Value = 1.825
$ws.Cells.Item(177, 21).Value = 2.5
$ws.Cells.Item(177, 22).Value = 1.75
$ws.Cells.Item(177, 23).Value = 1.95
$ws.Cells.Item(177, 24).Value = 0.55
$ws.Cells.Item(177, 25).Value = -1
$ws.Cells.Item(177, 26).Value = -1
$ws.Cells.Item(177, 27).Value = 0
$ws.Cells.Item(177, 28).Value = 0
$ws.Cells.Item(177, 29).Value = -1
$ws.Cells.Item(177, 30).Value = 0.95

# Row 178 <- data formerly on row 177
$ws.Cells.Item(178, 2).Value = 7217242
$ws.Cells.Item(178, 3).Value = 'Azerbaijan Premier League'
$ws.Cells.Item(178, 4).Value = 45437.5
$ws.Cells.Item(178, 5).Value = 'Zira IK'
$ws.Cells.Item(178, 6).Value = 'Sabail FC'
$ws.Cells.Item(178, 7).Value = 1
$ws.Cells.Item(178, 8).Value = 0
$ws.Cells.Item(178, 9).Value = 0
$ws.Cells.Item(178, 10).Value = 0
$ws.Cells.Item(178, 11).Value = 'H'
$ws.Cells.Item(178, 12).Value = 1.727
$ws.Cells.Item(178, 13).Value = 3.25
$ws.Cells.Item(178, 14).Value = 4.333
$ws.Cells.Item(178, 15).Value = 1.45
$ws.Cells.Item(178, 16).Value = 3.4
$ws.Cells.Item(178, 17).Value = 7
$ws.Cells.Item(178, 18).Value = -1
$ws.Cells.Item(178, 19).Value = 1.8
$ws.Cells.Item(178, 20).Value = 2
$ws.Cells.Item(178, 21).Value = 2.5
$ws.Cells.Item(178, 22).Value = 1.95
$ws.Cells.Item(178, 23).Value = 1.75
$ws.Cells.Item(178, 24).Value = 0.45
$ws.Cells.Item(178, 25).Value = -1
$ws.Cells.Item(178, 26).Value = -1
$ws.Cells.Item(178, 27).Value = 0
$ws.Cells.Item(178, 28).Value = 0
$ws.Cells.Item(178, 29).Value = -1
$ws.Cells.Item(178, 30).Value = 0.75
